$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers - activity types
$ws.Range("C1").Value = "ACTIVITY_TYPE_2"
$ws.Range("E1").Value = "ACTIVITY_TYPE_3"

# Row 3 - first sample question (food)
$ws.Range("A3").Value = "Is there Food?"
$ws.Range("B3").Value = "binary"

# Row 4 - second sample question (food quantity)
$ws.Range("A4").Value = "How much food is there?"
$ws.Range("B4").Value = "numeric"

# Row 2 - question / question type column headers for each activity type
$ws.Range("A2").Value = "QUESTION"
$ws.Range("B2").Value = "QUESTION_TYPE"
$ws.Range("C2").Value = "QUESTION"
$ws.Range("D2").Value = "QUESTION_TYPE"
$ws.Range("E2").Value = "QUESTION"
$ws.Range("F2").Value = "QUESTION_TYPE"

# Column width for new column F (closest achievable value to 27.1640625
# given the engine's 1/6-character ColumnWidth quantization)
$ws.Columns.Item(6).ColumnWidth = 26.25

# Selection matching target state
$ws.Range("E2:F2").Select()
